$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 corresponds to Start Duration "08:48:00" / End Duration "08:59:00"
# Update the Color and Color Code columns (C12, D12) from #673262 to #1e0e16
$ws.Range("C12").Value = "#1e0e16"
$ws.Range("D12").Value = "#1e0e16"

# Update the Last Changed On column (F12) from 12/03/2020 19:31:22 to 17/03/2020 00:05:33
$ws.Range("F12").Value = "17/03/2020 00:05:33"
